$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A1 holds a date serial; bump it by one day (45310 -> 45311)
$ws.Range("A1").Value = 45311

# Update unit prices in column D for rows 33-36
$ws.Range("D33").Value = 229.273
$ws.Range("D34").Value = 300
$ws.Range("D35").Value = 280.505
$ws.Range("D36").Value = 342.12
